$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column to the right of the existing "x" column (E),
# then move the "x" data into it so the new "cs_relative_to" column can
# take over the old E slot without disturbing its number formatting.
$ws.Columns.Item(6).Insert()

for ($r = 1; $r -le 21; $r++) {
    $srcCell = $ws.Cells.Item($r, 5)
    $dstCell = $ws.Cells.Item($r, 6)
    $dstCell.Value = $srcCell.Value2
}

$ws.Cells.Item(1, 5).Value = "cs_relative_to"
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "CS_1"
}

$ws.Columns.Item(5).ColumnWidth = 18.92

$ws.Range("G10").Select()
